# Apply weekly fruit/hortaliza data rotation across rows 2-6.
# The rows' date/variety/quality/volume/price/origin values get cyclically
# permuted: new row r gets the old values from row mapping[r].

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original values (columns D, K, L, M, N, O, P, R, S) for rows 2-6
# Note: use Value2 (not Value) -- this runtime mis-resolves the bare
# ".Value" property getter, so Value2 is used for reliable read/write.
$cols = @("D", "K", "L", "M", "N", "O", "P", "R", "S")
$orig = @{}
for ($r = 2; $r -le 6; $r++) {
    $orig[$r] = @{}
    foreach ($col in $cols) {
        $orig[$r][$col] = $ws.Range("$col$r").Value2
    }
}

# Mapping: new row -> source row (cyclic permutation observed in the diff)
$mapping = @{ 2 = 4; 3 = 6; 4 = 5; 5 = 3; 6 = 2 }

foreach ($r in $mapping.Keys) {
    $src = $mapping[$r]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value2 = $orig[$src][$col]
    }
}

$wb.Save()
